$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Q-calculation values for B2:I25 (24 rows x 8 columns)
$data = @(
    @(0.1980432023245938, 0.05992570125581027, 0.5570935873961507, 1.154973283651776, 0.2457951797412334, 0.007097295065019031, 0.4178298331263657, 0.2091803144137924),
    @(0.129319614615865, 0.1327175235045536, 0.5996073502723146, 1.156982352726067, 0.2296082113602821, 0.007100773653314739, 0.4072514853058511, 0.2118287548144155),
    @(0.1049460458845251, 0.1883645812486389, 0.6289481156376164, 1.153686981425362, 0.2167346828610983, 0.007103338927884559, 0.3994545999854755, 0.2088751067174635),
    @(0.1071165450562561, 0.2124430906695413, 0.6385821348136478, 1.145847968877939, 0.2099137122861406, 0.007104951343000625, 0.3945502658397297, 0.2011613001046939),
    @(0.1099451778327485, 0.2208187881726418, 0.6405144720547478, 1.140119102907062, 0.2069382590962053, 0.00710578882479528, 0.3920035605179626, 0.1954686543882198),
    @(0.1099451778327485, 0.2208187881726418, 0.6405144720547478, 1.140119102907062, 0.2069382590962053, 0.00710578882479528, 0.3920035605179626, 0.1954686543882198),
    @(0.1083127537332402, 0.2164418498585563, 0.6375969374377218, 1.139082022727119, 0.2076687929030149, 0.007105753988735217, 0.3921120514696886, 0.1944078254617973),
    @(0.1080509268275909, 0.1663755446811714, 0.6091431392423637, 1.137726679593879, 0.2182538029555467, 0.00710412065227192, 0.3970788789026836, 0.1927765353954893),
    @(0.2009211820402159, 0.05548766013464559, 0.5491326708474541, 1.143945872047216, 0.2445765414252417, 0.007098174303265969, 0.415161013353983, 0.1981256957906934),
    @(0.325494892923876, 0.06994868379532662, 0.4919787643866652, 1.148831499010274, 0.2698975923893778, 0.00709086700129691, 0.4373884903352633, 0.2018601469471601),
    @(0.385497614946199, 0.1193238713513806, 0.4621052098962299, 1.147859507827661, 0.28197941523589, 0.007087002215206147, 0.4483281121196138, 0.2002206033740415),
    @(0.4026780096839651, 0.1279734045038339, 0.438839036650387, 1.146565633982055, 0.2913509289758929, 0.007084285521503066, 0.452662658081405, 0.1983925647331062),
    @(0.3833331395007024, 0.1061330917975463, 0.4326354723103653, 1.144026943832149, 0.2929734307237618, 0.007084231194115395, 0.4516083648093541, 0.1957622157117314),
    @(0.3623958886899876, 0.087055496357616, 0.4342605740307829, 1.140000113074072, 0.2905880328008245, 0.007085337085007035, 0.4484228773970422, 0.1918557990537156),
    @(0.3546727873265608, 0.07934229807753763, 0.4379284626419198, 1.143785024681807, 0.2906148645058091, 0.007085237291593246, 0.4488586579819574, 0.1956738137861314),
    @(0.3477235793554523, 0.0723629051579892, 0.4424235872797243, 1.148708689506001, 0.2907663819322793, 0.007085035163405785, 0.4495941282307989, 0.2006309080446002),
    @(0.3113933629533152, 0.04530973521521984, 0.4595608210640599, 1.147796988848132, 0.2830879469910402, 0.007087516547276946, 0.4420797269358211, 0.20015221914545),
    @(0.2953674937055057, 0.03633048687199231, 0.4727244592829257, 1.149744531301193, 0.2783448924085293, 0.007088827378133221, 0.4378673174281558, 0.2023705197569689),
    @(0.3332898627407038, 0.06969052793151102, 0.472159171825476, 1.149334283070784, 0.2783922654332754, 0.007088531641788803, 0.4399455920120803, 0.2019352018892717),
    @(0.4307314280695401, 0.1576366558169681, 0.4458431718852264, 1.15257582479579, 0.2908423410002329, 0.00708366707309406, 0.4587906844567319, 0.2044436815672131),
    @(0.4932077061030367, 0.2102308743044842, 0.4239128545273678, 1.158973313412748, 0.3029826776564247, 0.007078792470677747, 0.4740813165371127, 0.2101026045794573),
    @(0.5375408583047505, 0.246639430384517, 0.4091454730240757, 1.165676481016805, 0.3123193893598374, 0.0070748209145583, 0.4861666131193653, 0.2162065739262563),
    @(0.5522093564310171, 0.2582737059387848, 0.4056748415957368, 1.170187418841742, 0.3158024755751556, 0.007073249915157342, 0.4909510011884161, 0.2205008683162257),
    @(0.4522825783961494, 0.1731792005658077, 0.4490949087987082, 1.172070928394247, 0.2976657088218361, 0.007098174303265969, 0.4489733892285723, 0.2236643829870656)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j  # Column B is index 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

Write-Output "Updated B2:I25 with new Q calculation values"